$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Simple odds updates on existing rows 3-21 (no structural change)
# ---------------------------------------------------------------
$ws.Range("G3").Value = 2.92
$ws.Range("J3").Value = 2.44
$ws.Range("R3").Value = 1.18
$ws.Range("W3").Value = 1.52

$ws.Range("P4").Value = 1.61

$ws.Range("F5").Value = 3.5
$ws.Range("G5").Value = 4.1
$ws.Range("H5").Value = 2.14
$ws.Range("I5").Value = 2.4
$ws.Range("K5").Value = 3.7
$ws.Range("N5").Value = 2.64
$ws.Range("P5").Value = 1.68
$ws.Range("Q5").Value = 2.02
$ws.Range("R5").Value = 1.22

$ws.Range("G6").Value = 4.1

$ws.Range("F8").Value = 1.04
$ws.Range("H8").Value = 1.04
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 1.01
$ws.Range("P8").Value = 1.76
$ws.Range("Q8").Value = 1.79

$ws.Range("K9").Value = 3.25

$ws.Range("H10").Value = 2.92
$ws.Range("K10").Value = 4.5

$ws.Range("G11").Value = 2.74
$ws.Range("H11").Value = 2.58
$ws.Range("J11").Value = 3.35
$ws.Range("K11").Value = 980

$ws.Range("I17").Value = 3.3
$ws.Range("J17").Value = 2.92
$ws.Range("Q17").Value = 2.54

$ws.Range("F19").Value = 2.14
$ws.Range("Q19").Value = 2.42
$ws.Range("T19").Value = 2.06
$ws.Range("AD19").Value = 22
$ws.Range("AH19").Value = 26
$ws.Range("AK19").Value = 32

$ws.Range("F20").Value = 1.81
$ws.Range("H20").Value = 2.62
$ws.Range("I20").Value = 5.3
$ws.Range("J20").Value = 2.2
$ws.Range("K20").Value = 3.85

$ws.Range("F21").Value = 1.81
$ws.Range("G21").Value = 1.83
$ws.Range("J21").Value = 4
$ws.Range("R21").Value = 1.46
$ws.Range("T21").Value = 1.75
$ws.Range("U21").Value = 2.22
$ws.Range("X21").Value = 19
$ws.Range("AG21").Value = 9.4

# ---------------------------------------------------------------
# 2) Insert a brand-new match row at row 22 (Spanish La Liga),
#    pushing the former rows 22-24 down to 23-25.
# ---------------------------------------------------------------
$ws.Rows(22).Insert()

function Set-TextCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "A22" "Spanish La Liga"
Set-TextCell "B22" "2025-11-24"
Set-TextCell "C22" "17:00:00"
Set-TextCell "D22" "Espanyol"
Set-TextCell "E22" "Sevilla"

$ws.Range("F22").Value = 2.12
$ws.Range("G22").Value = 2.14
$ws.Range("H22").Value = 4.1
$ws.Range("I22").Value = 4.2
$ws.Range("J22").Value = 3.45
$ws.Range("K22").Value = 3.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 1.1
$ws.Range("N22").Value = 3.25
$ws.Range("O22").Value = 1.43
$ws.Range("P22").Value = 1.78
$ws.Range("Q22").Value = 2.26
$ws.Range("R22").Value = 1.29
$ws.Range("S22").Value = 4.2
$ws.Range("T22").Value = 1.97
$ws.Range("U22").Value = 1.97
$ws.Range("V22").Value = 0
$ws.Range("W22").Value = 0
$ws.Range("X22").Value = 11
$ws.Range("Y22").Value = 13.5
$ws.Range("Z22").Value = 32
$ws.Range("AA22").Value = 120
$ws.Range("AB22").Value = 8.2
$ws.Range("AC22").Value = 7.6
$ws.Range("AD22").Value = 18
$ws.Range("AE22").Value = 1000
$ws.Range("AF22").Value = 12.5
$ws.Range("AG22").Value = 11
$ws.Range("AH22").Value = 21
$ws.Range("AI22").Value = 85
$ws.Range("AJ22").Value = 27
$ws.Range("AK22").Value = 25
$ws.Range("AL22").Value = 48
$ws.Range("AM22").Value = 160
$ws.Range("AN22").Value = 19.5
$ws.Range("AO22").Value = 1000

# ---------------------------------------------------------------
# 3) Tweaks on the rows that shifted down because of the insert
#    Row 23 = former row 22 (Brazilian Serie A)
#    Row 24 = former row 23 (Argentinian Primera Division)
#    Row 25 = former row 24 (Chilean Primera Division)
# ---------------------------------------------------------------
$ws.Range("J23").Value = 3.7

$ws.Range("F24").Value = 2.84

$ws.Range("G25").Value = 2.5
$ws.Range("H25").Value = 3.25
$ws.Range("I25").Value = 4.6
$ws.Range("J25").Value = 3.2
$ws.Range("K25").Value = 5.6
$ws.Range("N25").Value = 2.68
$ws.Range("P25").Value = 1.75
$ws.Range("Q25").Value = 1.87
$ws.Range("R25").Value = 1.29
$ws.Range("W25").Value = 1.66

$wb.Save()
